$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A - this shifts the existing
# A:D columns (segment name, RawActivations, PercActivations, totalActivation)
# one column to the right, becoming B:E. Excel moves each cell's own
# formatting along with it, so B:E keep the styles A:D used to have.
$ws.Range("A:A").Insert()

# New header for the inserted column, using the same header
# style (bold / centered / bordered) as the other header cells.
$ws.Range("B1").Value = "segments"
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)   # xlPasteFormats

# The old column A (now column B) held the segment name together with
# that same header style applied to every data row. After the shift the
# header style should no longer apply to the segment-name column, so
# clear formatting back to the default for the data rows.
$ws.Range("B2:B20").ClearFormats()

# Populate the new column A with a 0-based row index and give those
# cells the header style that column B's cells used to carry.
$names = @(
    "background",
    "back_bumper",
    "back_glass",
    "back_left_door",
    "back_left_light",
    "back_right_door",
    "back_right_light",
    "front_bumper",
    "front_glass",
    "front_left_door",
    "front_left_light",
    "front_right_door",
    "front_right_light",
    "hood",
    "left_mirror",
    "right_mirror",
    "tailgate",
    "trunk",
    "wheel"
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $names[$i]
}

$ws.Range("B1").Copy()
$ws.Range("A2:A20").PasteSpecial(-4122)   # xlPasteFormats
